$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: A3 was stored as text "20250308122656"; it becomes a real number.
$ws.Range("A3").Value = 20250308122656

# Row 4 (new)
$ws.Range("A4").Value = 20250308122815
$ws.Range("B4").Value = "rajas mhatre"
$ws.Range("C4").Value = 8548784834
$ws.Range("D4").Value = "jhsjbfhfbjshd sfdfd"

# Row 5 (new)
$ws.Range("A5").Value = 20250308123025
$ws.Range("B5").Value = "s pathak"
$ws.Range("C5").Value = 1232434
$ws.Range("D5").Value = "puneeee"

# Row 6 (new) - B/C/D are present but blank ("" inline strings)
$ws.Range("A6").Value = 20250308124203
$ws.Range("B6").Value = "'"
$ws.Range("B6").ClearFormats()
$ws.Range("C6").Value = "'"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").Value = "'"
$ws.Range("D6").ClearFormats()

# Row 7 (new) - B/C/D are present but blank ("" inline strings)
$ws.Range("A7").Value = 20250308131225
$ws.Range("B7").Value = "'"
$ws.Range("B7").ClearFormats()
$ws.Range("C7").Value = "'"
$ws.Range("C7").ClearFormats()
$ws.Range("D7").Value = "'"
$ws.Range("D7").ClearFormats()

# Row 8 (new) - B/C/D are present but blank ("" inline strings)
$ws.Range("A8").Value = 20250308131244
$ws.Range("B8").Value = "'"
$ws.Range("B8").ClearFormats()
$ws.Range("C8").Value = "'"
$ws.Range("C8").ClearFormats()
$ws.Range("D8").Value = "'"
$ws.Range("D8").ClearFormats()

# Row 9 (new) - B/C/D are present but blank ("" inline strings)
$ws.Range("A9").Value = 20250308131304
$ws.Range("B9").Value = "'"
$ws.Range("B9").ClearFormats()
$ws.Range("C9").Value = "'"
$ws.Range("C9").ClearFormats()
$ws.Range("D9").Value = "'"
$ws.Range("D9").ClearFormats()

# Row 10 (new) - B has text, C/D present but blank
$ws.Range("A10").Value = 20250308132522
$ws.Range("B10").Value = "rajas"
$ws.Range("C10").Value = "'"
$ws.Range("C10").ClearFormats()
$ws.Range("D10").Value = "'"
$ws.Range("D10").ClearFormats()

# Row 11 (new) - A11 and C11 stay text (numeric-looking strings), D11 blank
$ws.Range("A11").Value = "'20250308133046"
$ws.Range("A11").ClearFormats()
$ws.Range("B11").Value = "asas"
$ws.Range("C11").Value = "'3434"
$ws.Range("C11").ClearFormats()
$ws.Range("D11").Value = "'"
$ws.Range("D11").ClearFormats()
